$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 7).Value = 46.70430066666666
$ws.Cells.Item(2, 8).Value = 140.112902
$ws.Cells.Item(2, 9).Value = 0.4277960227396158
$ws.Cells.Item(2, 10).Value = 0.4350095176968582
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 13).Value = 33.211442
$ws.Cells.Item(2, 14).Value = 99.634326
$ws.Cells.Item(2, 15).Value = 0.211580186305583
$ws.Cells.Item(2, 16).Value = 0.2175281749633597
$ws.Cells.Item(2, 17).Value = 1551.117172741561
$ws.Cells.Item(2, 18).Value = 13960.05455467405
$ws.Cells.Item(2, 19).Value = 0.09051316219203531
$ws.Cells.Item(2, 20).Value = 0.0946268264762889

$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 7).Value = 46.70430066666666
$ws.Cells.Item(3, 8).Value = 140.112902
$ws.Cells.Item(3, 9).Value = 0.4277960227396158
$ws.Cells.Item(3, 10).Value = 0.4350095176968582
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 13).Value = 51.17424933333334
$ws.Cells.Item(3, 14).Value = 153.522748
$ws.Cells.Item(3, 15).Value = 0.3260158715178649
$ws.Cells.Item(3, 16).Value = 0.3351809012869699
$ws.Cells.Item(3, 17).Value = 2390.057527254966
$ws.Cells.Item(3, 18).Value = 21510.5177452947
$ws.Cells.Item(3, 19).Value = 0.1394682931853322
$ws.Cells.Item(3, 20).Value = 0.145806882210043

$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 7).Value = 46.70430066666666
$ws.Cells.Item(4, 8).Value = 140.112902
$ws.Cells.Item(4, 9).Value = 0.4277960227396158
$ws.Cells.Item(4, 10).Value = 0.4350095176968582
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 13).Value = 22.19450366666666
$ws.Cells.Item(4, 14).Value = 66.58351099999999
$ws.Cells.Item(4, 15).Value = 0.1413945597650736
$ws.Cells.Item(4, 16).Value = 0.1453694746776606
$ws.Cells.Item(4, 17).Value = 1036.578772395435
$ws.Cells.Item(4, 18).Value = 9329.20895155892
$ws.Cells.Item(4, 19).Value = 0.0604880303045174
$ws.Cells.Item(4, 20).Value = 0.06323710506737479

$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 7).Value = 46.70430066666666
$ws.Cells.Item(5, 8).Value = 140.112902
$ws.Cells.Item(5, 9).Value = 0.4277960227396158
$ws.Cells.Item(5, 10).Value = 0.4350095176968582
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 13).Value = 37.51216133333333
$ws.Cells.Item(5, 14).Value = 112.536484
$ws.Cells.Item(5, 15).Value = 0.2389787857941174
$ws.Cells.Item(5, 16).Value = 0.2456970098971044
$ws.Cells.Item(5, 17).Value = 1751.979261568507
$ws.Cells.Item(5, 18).Value = 15767.81335411657
$ws.Cells.Item(5, 19).Value = 0.102234174081866
$ws.Cells.Item(5, 20).Value = 0.1068805377748996

$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 7).Value = 46.70430066666666
$ws.Cells.Item(6, 8).Value = 140.112902
$ws.Cells.Item(6, 9).Value = 0.4277960227396158
$ws.Cells.Item(6, 10).Value = 0.4350095176968582
$ws.Cells.Item(6, 11).Value = 2.0
$ws.Cells.Item(6, 13).Value = 12.8762265
$ws.Cells.Item(6, 14).Value = 25.752453
$ws.Cells.Item(6, 15).Value = 0.08203059661736112
$ws.Cells.Item(6, 16).Value = 0.05622443917490542
$ws.Cells.Item(6, 17).Value = 601.375153908101
$ws.Cells.Item(6, 18).Value = 3608.250923448606
$ws.Cells.Item(6, 19).Value = 0.03509236297586486
$ws.Cells.Item(6, 20).Value = 0.02445816616825195

$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 7).Value = 14.05160533333333
$ws.Cells.Item(7, 8).Value = 42.154816
$ws.Cells.Item(7, 9).Value = 0.1287080801746603
$ws.Cells.Item(7, 10).Value = 0.1308783553477452
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 13).Value = 33.211442
$ws.Cells.Item(7, 14).Value = 99.634326
$ws.Cells.Item(7, 15).Value = 0.211580186305583
$ws.Cells.Item(7, 16).Value = 0.2175281749633597
$ws.Cells.Item(7, 17).Value = 466.6740755348906
$ws.Cells.Item(7, 18).Value = 4200.066679814016
$ws.Cells.Item(7, 19).Value = 0.02723207958238854
$ws.Cells.Item(7, 20).Value = 0.02846972978100109

$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 7).Value = 14.05160533333333
$ws.Cells.Item(8, 8).Value = 42.154816
$ws.Cells.Item(8, 9).Value = 0.1287080801746603
$ws.Cells.Item(8, 10).Value = 0.1308783553477452
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 13).Value = 51.17424933333334
$ws.Cells.Item(8, 14).Value = 153.522748
$ws.Cells.Item(8, 15).Value = 0.3260158715178649
$ws.Cells.Item(8, 16).Value = 0.3351809012869699
$ws.Cells.Item(8, 17).Value = 719.0803548615964
$ws.Cells.Item(8, 18).Value = 6471.723193754367
$ws.Cells.Item(8, 19).Value = 0.0419608769295331
$ws.Cells.Item(8, 20).Value = 0.04386792510441356

$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 7).Value = 14.05160533333333
$ws.Cells.Item(9, 8).Value = 42.154816
$ws.Cells.Item(9, 9).Value = 0.1287080801746603
$ws.Cells.Item(9, 10).Value = 0.1308783553477452
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 13).Value = 22.19450366666666
$ws.Cells.Item(9, 14).Value = 66.58351099999999
$ws.Cells.Item(9, 15).Value = 0.1413945597650736
$ws.Cells.Item(9, 16).Value = 0.1453694746776606
$ws.Cells.Item(9, 17).Value = 311.8684060932195
$ws.Cells.Item(9, 18).Value = 2806.815654838975
$ws.Cells.Item(9, 19).Value = 0.01819862233450389
$ws.Cells.Item(9, 20).Value = 0.01902571776357791

$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 7).Value = 14.05160533333333
$ws.Cells.Item(10, 8).Value = 42.154816
$ws.Cells.Item(10, 9).Value = 0.1287080801746603
$ws.Cells.Item(10, 10).Value = 0.1308783553477452
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 13).Value = 37.51216133333333
$ws.Cells.Item(10, 14).Value = 112.536484
$ws.Cells.Item(10, 15).Value = 0.2389787857941174
$ws.Cells.Item(10, 16).Value = 0.2456970098971044
$ws.Cells.Item(10, 17).Value = 527.1060862563271
$ws.Cells.Item(10, 18).Value = 4743.954776306944
$ws.Cells.Item(10, 19).Value = 0.03075850072203222
$ws.Cells.Item(10, 20).Value = 0.03215642056919171

$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 7).Value = 14.05160533333333
$ws.Cells.Item(11, 8).Value = 42.154816
$ws.Cells.Item(11, 9).Value = 0.1287080801746603
$ws.Cells.Item(11, 10).Value = 0.1308783553477452
$ws.Cells.Item(11, 11).Value = 2.0
$ws.Cells.Item(11, 13).Value = 12.8762265
$ws.Cells.Item(11, 14).Value = 25.752453
$ws.Cells.Item(11, 15).Value = 0.08203059661736112
$ws.Cells.Item(11, 16).Value = 0.05622443917490542
$ws.Cells.Item(11, 17).Value = 180.931652960608
$ws.Cells.Item(11, 18).Value = 1085.589917763648
$ws.Cells.Item(11, 19).Value = 0.01055800060620253
$ws.Cells.Item(11, 20).Value = 0.007358562129560958

$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 7).Value = 15.248849
$ws.Cells.Item(12, 8).Value = 45.746547
$ws.Cells.Item(12, 9).Value = 0.1396744381232708
$ws.Cells.Item(12, 10).Value = 0.1420296279836289
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 13).Value = 33.211442
$ws.Cells.Item(12, 14).Value = 99.634326
$ws.Cells.Item(12, 15).Value = 0.211580186305583
$ws.Cells.Item(12, 16).Value = 0.2175281749633597
$ws.Cells.Item(12, 17).Value = 506.4362641302579
$ws.Cells.Item(12, 18).Value = 4557.926377172322
$ws.Cells.Item(12, 19).Value = 0.02955234364024925
$ws.Cells.Item(12, 20).Value = 0.03089544576600373

$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 7).Value = 15.248849
$ws.Cells.Item(13, 8).Value = 45.746547
$ws.Cells.Item(13, 9).Value = 0.1396744381232708
$ws.Cells.Item(13, 10).Value = 0.1420296279836289
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 13).Value = 51.17424933333334
$ws.Cells.Item(13, 14).Value = 153.522748
$ws.Cells.Item(13, 15).Value = 0.3260158715178649
$ws.Cells.Item(13, 16).Value = 0.3351809012869699
$ws.Cells.Item(13, 17).Value = 780.3484007723507
$ws.Cells.Item(13, 18).Value = 7023.135606951157
$ws.Cells.Item(13, 19).Value = 0.04553608367352622
$ws.Cells.Item(13, 20).Value = 0.0476056187170058

$ws.Cells.Item(14, 5).Value = 3.0
$ws.Cells.Item(14, 7).Value = 15.248849
$ws.Cells.Item(14, 8).Value = 45.746547
$ws.Cells.Item(14, 9).Value = 0.1396744381232708
$ws.Cells.Item(14, 10).Value = 0.1420296279836289
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 13).Value = 22.19450366666666
$ws.Cells.Item(14, 14).Value = 66.58351099999999
$ws.Cells.Item(14, 15).Value = 0.1413945597650736
$ws.Cells.Item(14, 16).Value = 0.1453694746776606
$ws.Cells.Item(14, 17).Value = 338.4406350429463
$ws.Cells.Item(14, 18).Value = 3045.965715386516
$ws.Cells.Item(14, 19).Value = 0.01974920568887389
$ws.Cells.Item(14, 20).Value = 0.0206467724086437

$ws.Cells.Item(15, 5).Value = 3.0
$ws.Cells.Item(15, 7).Value = 15.248849
$ws.Cells.Item(15, 8).Value = 45.746547
$ws.Cells.Item(15, 9).Value = 0.1396744381232708
$ws.Cells.Item(15, 10).Value = 0.1420296279836289
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 13).Value = 37.51216133333333
$ws.Cells.Item(15, 14).Value = 112.536484
$ws.Cells.Item(15, 15).Value = 0.2389787857941174
$ws.Cells.Item(15, 16).Value = 0.2456970098971044
$ws.Cells.Item(15, 17).Value = 572.0172838356386
$ws.Cells.Item(15, 18).Value = 5148.155554520748
$ws.Cells.Item(15, 19).Value = 0.03337922762917482
$ws.Cells.Item(15, 20).Value = 0.03489625491237574

$ws.Cells.Item(16, 5).Value = 3.0
$ws.Cells.Item(16, 7).Value = 15.248849
$ws.Cells.Item(16, 8).Value = 45.746547
$ws.Cells.Item(16, 9).Value = 0.1396744381232708
$ws.Cells.Item(16, 10).Value = 0.1420296279836289
$ws.Cells.Item(16, 11).Value = 2.0
$ws.Cells.Item(16, 13).Value = 12.8762265
$ws.Cells.Item(16, 14).Value = 25.752453
$ws.Cells.Item(16, 15).Value = 0.08203059661736112
$ws.Cells.Item(16, 16).Value = 0.05622443917490542
$ws.Cells.Item(16, 17).Value = 196.3476335882985
$ws.Cells.Item(16, 18).Value = 1178.085801529791
$ws.Cells.Item(16, 19).Value = 0.01145757749144659
$ws.Cells.Item(16, 20).Value = 0.00798553617959999

$ws.Cells.Item(17, 5).Value = 3.0
$ws.Cells.Item(17, 7).Value = 27.73836633333333
$ws.Cells.Item(17, 8).Value = 83.215099
$ws.Cells.Item(17, 9).Value = 0.2540743063339262
$ws.Cells.Item(17, 10).Value = 0.258358506350017
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 13).Value = 33.211442
$ws.Cells.Item(17, 14).Value = 99.634326
$ws.Cells.Item(17, 15).Value = 0.211580186305583
$ws.Cells.Item(17, 16).Value = 0.2175281749633597
$ws.Cells.Item(17, 17).Value = 921.2311446542526
$ws.Cells.Item(17, 18).Value = 8291.080301888274
$ws.Cells.Item(17, 19).Value = 0.05375708906959387
$ws.Cells.Item(17, 20).Value = 0.05620025437257877

$ws.Cells.Item(18, 5).Value = 3.0
$ws.Cells.Item(18, 7).Value = 27.73836633333333
$ws.Cells.Item(18, 8).Value = 83.215099
$ws.Cells.Item(18, 9).Value = 0.2540743063339262
$ws.Cells.Item(18, 10).Value = 0.258358506350017
$ws.Cells.Item(18, 11).Value = 3.0
$ws.Cells.Item(18, 13).Value = 51.17424933333334
$ws.Cells.Item(18, 14).Value = 153.522748
$ws.Cells.Item(18, 15).Value = 0.3260158715178649
$ws.Cells.Item(18, 16).Value = 0.3351809012869699
$ws.Cells.Item(18, 17).Value = 1419.490074841339
$ws.Cells.Item(18, 18).Value = 12775.41067357205
$ws.Cells.Item(18, 19).Value = 0.08283225640975193
$ws.Cells.Item(18, 20).Value = 0.08659683701355404

$ws.Cells.Item(19, 5).Value = 3.0
$ws.Cells.Item(19, 7).Value = 27.73836633333333
$ws.Cells.Item(19, 8).Value = 83.215099
$ws.Cells.Item(19, 9).Value = 0.2540743063339262
$ws.Cells.Item(19, 10).Value = 0.258358506350017
$ws.Cells.Item(19, 11).Value = 3.0
$ws.Cells.Item(19, 13).Value = 22.19450366666666
$ws.Cells.Item(19, 14).Value = 66.58351099999999
$ws.Cells.Item(19, 15).Value = 0.1413945597650736
$ws.Cells.Item(19, 16).Value = 0.1453694746776606
$ws.Cells.Item(19, 17).Value = 615.6392732925098
$ws.Cells.Item(19, 18).Value = 5540.753459632588
$ws.Cells.Item(19, 19).Value = 0.03592472469170196
$ws.Cells.Item(19, 20).Value = 0.03755744034660701

$ws.Cells.Item(20, 5).Value = 3.0
$ws.Cells.Item(20, 7).Value = 27.73836633333333
$ws.Cells.Item(20, 8).Value = 83.215099
$ws.Cells.Item(20, 9).Value = 0.2540743063339262
$ws.Cells.Item(20, 10).Value = 0.258358506350017
$ws.Cells.Item(20, 11).Value = 3.0
$ws.Cells.Item(20, 13).Value = 37.51216133333333
$ws.Cells.Item(20, 14).Value = 112.536484
$ws.Cells.Item(20, 15).Value = 0.2389787857941174
$ws.Cells.Item(20, 16).Value = 0.2456970098971044
$ws.Cells.Item(20, 17).Value = 1040.526073019102
$ws.Cells.Item(20, 18).Value = 9364.734657171915
$ws.Cells.Item(20, 19).Value = 0.0607183692291643
$ws.Cells.Item(20, 20).Value = 0.06347791249168123

$ws.Cells.Item(21, 5).Value = 3.0
$ws.Cells.Item(21, 7).Value = 27.73836633333333
$ws.Cells.Item(21, 8).Value = 83.215099
$ws.Cells.Item(21, 9).Value = 0.2540743063339262
$ws.Cells.Item(21, 10).Value = 0.258358506350017
$ws.Cells.Item(21, 11).Value = 2.0
$ws.Cells.Item(21, 13).Value = 12.8762265
$ws.Cells.Item(21, 14).Value = 25.752453
$ws.Cells.Item(21, 15).Value = 0.08203059661736112
$ws.Cells.Item(21, 16).Value = 0.05622443917490542
$ws.Cells.Item(21, 17).Value = 357.1654876479745
$ws.Cells.Item(21, 18).Value = 2142.992925887847
$ws.Cells.Item(21, 19).Value = 0.02084186693371414
$ws.Cells.Item(21, 20).Value = 0.01452606212559595

$ws.Cells.Item(22, 5).Value = 2.0
$ws.Cells.Item(22, 7).Value = 5.431107000000001
$ws.Cells.Item(22, 8).Value = 10.862214
$ws.Cells.Item(22, 9).Value = 0.0497471526285271
$ws.Cells.Item(22, 10).Value = 0.03372399262175058
$ws.Cells.Item(22, 11).Value = 3.0
$ws.Cells.Item(22, 13).Value = 33.211442
$ws.Cells.Item(22, 14).Value = 99.634326
$ws.Cells.Item(22, 15).Value = 0.211580186305583
$ws.Cells.Item(22, 16).Value = 0.2175281749633597
$ws.Cells.Item(22, 17).Value = 180.374895126294
$ws.Cells.Item(22, 18).Value = 1082.249370757764
$ws.Cells.Item(22, 19).Value = 0.01052551182131604
$ws.Cells.Item(22, 20).Value = 0.007335918567487212

$ws.Cells.Item(23, 5).Value = 2.0
$ws.Cells.Item(23, 7).Value = 5.431107000000001
$ws.Cells.Item(23, 8).Value = 10.862214
$ws.Cells.Item(23, 9).Value = 0.0497471526285271
$ws.Cells.Item(23, 10).Value = 0.03372399262175058
$ws.Cells.Item(23, 11).Value = 3.0
$ws.Cells.Item(23, 13).Value = 51.17424933333334
$ws.Cells.Item(23, 14).Value = 153.522748
$ws.Cells.Item(23, 15).Value = 0.3260158715178649
$ws.Cells.Item(23, 16).Value = 0.3351809012869699
$ws.Cells.Item(23, 17).Value = 277.932823774012
$ws.Cells.Item(23, 18).Value = 1667.596942644072
$ws.Cells.Item(23, 19).Value = 0.01621836131972151
$ws.Cells.Item(23, 20).Value = 0.01130363824195348

$ws.Cells.Item(24, 5).Value = 2.0
$ws.Cells.Item(24, 7).Value = 5.431107000000001
$ws.Cells.Item(24, 8).Value = 10.862214
$ws.Cells.Item(24, 9).Value = 0.0497471526285271
$ws.Cells.Item(24, 10).Value = 0.03372399262175058
$ws.Cells.Item(24, 11).Value = 3.0
$ws.Cells.Item(24, 13).Value = 22.19450366666666
$ws.Cells.Item(24, 14).Value = 66.58351099999999
$ws.Cells.Item(24, 15).Value = 0.1413945597650736
$ws.Cells.Item(24, 16).Value = 0.1453694746776606
$ws.Cells.Item(24, 17).Value = 120.540724225559
$ws.Cells.Item(24, 18).Value = 723.244345353354
$ws.Cells.Item(24, 19).Value = 0.007033976745476516
$ws.Cells.Item(24, 20).Value = 0.004902439091457184

$ws.Cells.Item(25, 5).Value = 2.0
$ws.Cells.Item(25, 7).Value = 5.431107000000001
$ws.Cells.Item(25, 8).Value = 10.862214
$ws.Cells.Item(25, 9).Value = 0.0497471526285271
$ws.Cells.Item(25, 10).Value = 0.03372399262175058
$ws.Cells.Item(25, 11).Value = 3.0
$ws.Cells.Item(25, 13).Value = 37.51216133333333
$ws.Cells.Item(25, 14).Value = 112.536484
$ws.Cells.Item(25, 15).Value = 0.2389787857941174
$ws.Cells.Item(25, 16).Value = 0.2456970098971044
$ws.Cells.Item(25, 17).Value = 203.732562002596
$ws.Cells.Item(25, 18).Value = 1222.395372015576
$ws.Cells.Item(25, 19).Value = 0.01188851413188004
$ws.Cells.Item(25, 20).Value = 0.008285884148956129

$ws.Cells.Item(26, 5).Value = 2.0
$ws.Cells.Item(26, 7).Value = 5.431107000000001
$ws.Cells.Item(26, 8).Value = 10.862214
$ws.Cells.Item(26, 9).Value = 0.0497471526285271
$ws.Cells.Item(26, 10).Value = 0.03372399262175058
$ws.Cells.Item(26, 11).Value = 2.0
$ws.Cells.Item(26, 13).Value = 12.8762265
$ws.Cells.Item(26, 14).Value = 25.752453
$ws.Cells.Item(26, 15).Value = 0.08203059661736112
$ws.Cells.Item(26, 16).Value = 0.05622443917490542
$ws.Cells.Item(26, 17).Value = 69.93216387773552
$ws.Cells.Item(26, 18).Value = 279.7286555109421
$ws.Cells.Item(26, 19).Value = 0.004080788610133002
$ws.Cells.Item(26, 20).Value = 0.001896112571896575
